# hr_boot: progress button and adding data
# The "Creative City" sheet contains a couple of instructional placeholder
# strings ("keep it empty " / "keep it empty" / "filled automatically ")
# that should not actually hold text - clear them out so the cells are
# blank (keeping their existing formatting/style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Creative City")

$ws.Range("G3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("Q3").Value = ""

# Leave the cursor where the author last left it.
$ws.Range("F16").Select() | Out-Null
